$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap columns A (Product_Name) and B (Category) ---
# Header row
$ws.Range("A1").Value = "Product_Name"
$ws.Range("B1").Value = "Category"

# Data rows: new A = old B (product name), new B = old A (category)
$ws.Range("A2").Value = "Carhartt Men's Acrylic Watch Hat"
$ws.Range("B2").Value = "Accessories"

$ws.Range("A3").Value = "CRZ YOGA Butterluxe High Waisted Leggings"
$ws.Range("B3").Value = "Clothing"

$ws.Range("A4").Value = "PAVOI 14K Gold Plated Lightweight Hoops"
$ws.Range("B4").Value = "Accessories"

$ws.Range("A5").Value = "Gildan Adult Fleece Hooded Sweatshirt"
$ws.Range("B5").Value = "Clothing"

$ws.Range("A6").Value = "Crocs Unisex Classic Clogs"
$ws.Range("B6").Value = "Shoes"

# --- Column F: "Pinterest_Title" -> "Source" with new values ---
$ws.Range("F1").Value = "Source"
$ws.Range("F2").Value = "CC Picks"
$ws.Range("F3").Value = "Toronto Base"
$ws.Range("F4").Value = "Amazon Top Choice"
$ws.Range("F5").Value = "Amazon Top Choice"
$ws.Range("F6").Value = "Amazon Top Choice"

# --- Give F2:F3 a distinct (fill-flagged) style, matching the source workbook ---
$ws.Range("F2:F3").Interior.ColorIndex = -4142

# --- New (empty, but styled) row 9, cell C9 matching the existing body style ---
$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths: A/B swap, F narrows ---
$ws.Columns.Item(1).ColumnWidth = 23.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
$ws.Columns.Item(6).ColumnWidth = 16.666666666666668

# --- Selection ---
$ws.Range("D13").Select()
